$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '38.017.01'
$ws.Range("E2").Value = '  +2.61%  '

# Row 3
$ws.Range("D3").Value = '2.048.60'
$ws.Range("E3").Value = '  +1.68%  '

# Row 4
$ws.Range("E4").Value = '  +0.00%  '

# Row 5
$ws.Range("D5").Value = '''228.74'
$ws.Range("E5").Value = '  +0.66%  '

# Row 6
$ws.Range("E6").Value = '  +0.95%  '

# Row 7
$ws.Range("D7").Value = '''60.54'
$ws.Range("E7").Value = '  +8.57%  '

# Row 8
$ws.Range("E8").Value = '  -0.01%  '

# Row 9
$ws.Range("E9").Value = '  +2.40%  '

# Row 10
$ws.Range("D10").Value = '''0.0802'
$ws.Range("E10").Value = '  +3.07%  '

# Row 11
$ws.Range("E11").Value = '  +1.91%  '

# Row 12
$ws.Range("B12").Value = 'Chainlink'
$ws.Range("C12").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D12").Value = '''14.72'
$ws.Range("E12").Value = '  +3.56%  '

# Row 13
$ws.Range("B13").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C13").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D13").Value = '2.352.98'
$ws.Range("E13").Value = '  +1.50%  '

# Row 14
$ws.Range("D14").Value = '''21.03'
$ws.Range("E14").Value = '  +5.69%  '

# Row 15
$ws.Range("D15").Value = '''5.32'
$ws.Range("E15").Value = '  +3.44%  '

# Row 16
$ws.Range("D16").Value = '''0.756'
$ws.Range("E16").Value = '  +2.81%  '

# Row 17
$ws.Range("D17").Value = '2.054.17'
$ws.Range("E17").Value = '  +2.09%  '

# Row 18
$ws.Range("D18").Value = '37.972.35'
$ws.Range("E18").Value = '  +2.63%  '

# Row 19
$ws.Range("D19").Value = '''6.27'
$ws.Range("E19").Value = '  +2.48%  '

# Row 20
$ws.Range("E20").Value = '  +1.21%  '

# Row 21
$ws.Range("D21").Value = '0.0₃0830'
$ws.Range("E21").Value = '  +2.03%  '

# Row 22
$ws.Range("D22").Value = '''225.67'
$ws.Range("E22").Value = '  +1.37%  '

# Row 23
$ws.Range("E23").Value = '  +0.08%  '

# Row 24
$ws.Range("E24").Value = '  -0.32%  '

# Row 25
$ws.Range("E25").Value = '  +1.78%  '

# Row 26
$ws.Range("D26").Value = '''165.50'
$ws.Range("E26").Value = '  +1.39%  '

# Row 27
$ws.Range("D27").Value = '''9.18'
$ws.Range("E27").Value = '  +2.36%  '

# Row 28
$ws.Range("D28").Value = '''0.133'
$ws.Range("E28").Value = '  +4.76%  '

# Row 29
$ws.Range("D29").Value = '''19.00'
$ws.Range("E29").Value = '  +1.79%  '

# Row 30
$ws.Range("D30").Value = '''1.30'
$ws.Range("E30").Value = '  -0.12%  '

# Row 31
$ws.Range("D31").Value = '''0.120'
$ws.Range("E31").Value = '  +2.15%  '

# Row 32
$ws.Range("D32").Value = '''4.51'
$ws.Range("E32").Value = '  +1.95%  '

# Row 33
$ws.Range("D33").Value = '''4.56'
$ws.Range("E33").Value = '  +2.36%  '

# Row 34
$ws.Range("D34").Value = '''2.03'
$ws.Range("E34").Value = '  +8.62%  '

# Row 35
$ws.Range("D35").Value = '''0.0603'
$ws.Range("E35").Value = '  +0.30%  '

# Row 36
$ws.Range("D36").Value = '''6.25'
$ws.Range("E36").Value = '  +14.81%  '

# Row 37
$ws.Range("D37").Value = '''2.30'
$ws.Range("E37").Value = '  -1.86%  '

# Row 38
$ws.Range("D38").Value = '''3.27'
$ws.Range("E38").Value = '  +3.70%  '

# Row 39
$ws.Range("E39").Value = '  -0.13%  '

# Row 40
$ws.Range("D40").Value = '1.520.23'
$ws.Range("E40").Value = '  +3.54%  '

# Row 41
$ws.Range("D41").Value = '''97.27'
$ws.Range("E41").Value = '  +3.56%  '

# Row 42
$ws.Range("B42").Value = 'InjectiveProtocol'
$ws.Range("C42").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D42").Value = '''16.87'
$ws.Range("E42").Value = '  +4.22%  '

# Row 43
$ws.Range("B43").Value = 'VeChain'
$ws.Range("C43").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D43").Value = '''0.0215'
$ws.Range("E43").Value = '  +1.36%  '

# Row 44
$ws.Range("B44").Value = 'HuobiToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D44").Value = '''2.87'
$ws.Range("E44").Value = '  +3.51%  '

# Row 45
$ws.Range("D45").Value = '''0.0922'
$ws.Range("E45").Value = '  +1.25%  '

# Row 46
$ws.Range("E46").Value = '  +2.17%  '

# Row 47
$ws.Range("E47").Value = '  -4.99%  '

# Row 48
$ws.Range("E48").Value = '  +1.06%  '

# Row 50
$ws.Range("D50").Value = '''7.00'
$ws.Range("E50").Value = '  -0.50%  '

# Row 51
$ws.Range("D51").Value = '2.240.96'
$ws.Range("E51").Value = '  +1.68%  '
